# Generate Report for Handoff
# Adds two new localization-status rows (512be6a6-... and ae7272b2-...)
# to the Overview sheet and the per-locale (zh-cn / de-de) detail sheets,
# then grows the backing tables to cover the new rows.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------
# Sheet: Overview  (columns A:G, new rows 6 & 7)
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A6").Value = "512be6a6-b0d4-4eaf-af6c-cb2efeee8c34.md"
$ov.Range("C6").Value = ".md"
$ov.Range("E6").Value = "Ready for handoff"
$ov.Range("F6").Value = "Ready for handoff"
$ov.Range("G6").Value = "2016-08-29 16:45:40"

$ov.Range("A7").Value = "ae7272b2-561f-4c2a-a2c1-f61fde1279a2.md"
$ov.Range("C7").Value = ".md"
$ov.Range("E7").Value = "Ready for handoff"
$ov.Range("F7").Value = "Ready for handoff"
$ov.Range("G7").Value = "2016-08-29 16:45:40"

$ov.Range("G6:G7").NumberFormat = $dateFmt

# Path And Name column (B) carries the hyperlink, same as rows 2-5.
$ov.Range("B6").Value = "e2e\512be6a6-b0d4-4eaf-af6c-cb2efeee8c34.md"
$ov.Hyperlinks.Add($ov.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/512be6a6b0d44eafaf6ccb2efeee8c34512be6a/e2e/512be6a6-b0d4-4eaf-af6c-cb2efeee8c34.md", "", "", "e2e\512be6a6-b0d4-4eaf-af6c-cb2efeee8c34.md")

$ov.Range("B7").Value = "e2e\ae7272b2-561f-4c2a-a2c1-f61fde1279a2.md"
$ov.Hyperlinks.Add($ov.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ae7272b2561f4c2aa2c1f61fde1279a2ae7272b/e2e/ae7272b2-561f-4c2a-a2c1-f61fde1279a2.md", "", "", "e2e\ae7272b2-561f-4c2a-a2c1-f61fde1279a2.md")

$ovTbl = $ov.ListObjects.Item(1)
$ovTbl.Resize($ov.Range("A1:G7"))

# ---------------------------------------------------------------
# Sheet: zh-cn  (columns A:P, new rows 6 & 7)
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B6").Value = ".md"
$zh.Range("C6").Value = "Ready for handoff"
$zh.Range("D6").Value = "e2e"
$zh.Range("E6").Value = "ht"
$zh.Range("F6").Value = "False"
$zh.Range("G6").Value = "512be6a6-b0d4-4eaf-af6c-cb2efeee8c34.1e71e1b25cb548868f37319af2ff13012a90ca13.zh-cn.xlf"
$zh.Range("H6").Value = "2016-08-29 16:45:35"
$zh.Range("K6").Value = "0001-01-01 00:00:00"
$zh.Range("M6").Value = "True"
$zh.Range("O6").Value = "False"

$zh.Range("B7").Value = ".md"
$zh.Range("C7").Value = "Ready for handoff"
$zh.Range("D7").Value = "e2e"
$zh.Range("E7").Value = "ht"
$zh.Range("F7").Value = "False"
$zh.Range("G7").Value = "ae7272b2-561f-4c2a-a2c1-f61fde1279a2.f0144ba33e22fd596802a3b78095f849d7467bb5.zh-cn.xlf"
$zh.Range("H7").Value = "2016-08-29 16:45:35"
$zh.Range("K7").Value = "0001-01-01 00:00:00"
$zh.Range("M7").Value = "True"
$zh.Range("O7").Value = "False"

$zh.Range("H6:H7").NumberFormat = $dateFmt
$zh.Range("K6:K7").NumberFormat = $dateFmt

# Source File Name column (A) carries the hyperlink, same as rows 2-5.
$zh.Range("A6").Value = "512be6a6-b0d4-4eaf-af6c-cb2efeee8c34.md"
$zh.Hyperlinks.Add($zh.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/512be6a6b0d44eafaf6ccb2efeee8c34512be6a/e2e/512be6a6-b0d4-4eaf-af6c-cb2efeee8c34.md", "", "", "512be6a6-b0d4-4eaf-af6c-cb2efeee8c34.md")

$zh.Range("A7").Value = "ae7272b2-561f-4c2a-a2c1-f61fde1279a2.md"
$zh.Hyperlinks.Add($zh.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ae7272b2561f4c2aa2c1f61fde1279a2ae7272b/e2e/ae7272b2-561f-4c2a-a2c1-f61fde1279a2.md", "", "", "ae7272b2-561f-4c2a-a2c1-f61fde1279a2.md")

$zhTbl = $zh.ListObjects.Item(1)
$zhTbl.Resize($zh.Range("A1:P7"))

# ---------------------------------------------------------------
# Sheet: de-de  (columns A:P, new rows 6 & 7)
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B6").Value = ".md"
$de.Range("C6").Value = "Ready for handoff"
$de.Range("D6").Value = "e2e"
$de.Range("E6").Value = "ht"
$de.Range("F6").Value = "False"
$de.Range("G6").Value = "512be6a6-b0d4-4eaf-af6c-cb2efeee8c34.1e71e1b25cb548868f37319af2ff13012a90ca13.de-de.xlf"
$de.Range("H6").Value = "2016-08-29 16:45:40"
$de.Range("K6").Value = "0001-01-01 00:00:00"
$de.Range("M6").Value = "True"
$de.Range("O6").Value = "False"

$de.Range("B7").Value = ".md"
$de.Range("C7").Value = "Ready for handoff"
$de.Range("D7").Value = "e2e"
$de.Range("E7").Value = "ht"
$de.Range("F7").Value = "False"
$de.Range("G7").Value = "ae7272b2-561f-4c2a-a2c1-f61fde1279a2.f0144ba33e22fd596802a3b78095f849d7467bb5.de-de.xlf"
$de.Range("H7").Value = "2016-08-29 16:45:40"
$de.Range("K7").Value = "0001-01-01 00:00:00"
$de.Range("M7").Value = "True"
$de.Range("O7").Value = "False"

$de.Range("H6:H7").NumberFormat = $dateFmt
$de.Range("K6:K7").NumberFormat = $dateFmt

# Source File Name column (A) carries the hyperlink, same as rows 2-5.
$de.Range("A6").Value = "512be6a6-b0d4-4eaf-af6c-cb2efeee8c34.md"
$de.Hyperlinks.Add($de.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/512be6a6b0d44eafaf6ccb2efeee8c34512be6a/e2e/512be6a6-b0d4-4eaf-af6c-cb2efeee8c34.md", "", "", "512be6a6-b0d4-4eaf-af6c-cb2efeee8c34.md")

$de.Range("A7").Value = "ae7272b2-561f-4c2a-a2c1-f61fde1279a2.md"
$de.Hyperlinks.Add($de.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ae7272b2561f4c2aa2c1f61fde1279a2ae7272b/e2e/ae7272b2-561f-4c2a-a2c1-f61fde1279a2.md", "", "", "ae7272b2-561f-4c2a-a2c1-f61fde1279a2.md")

$deTbl = $de.ListObjects.Item(1)
$deTbl.Resize($de.Range("A1:P7"))
